$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2; existing rows 2-11 shift down to 3-12.
$ws.Rows(2).Insert()

# The freshly inserted row inherits the (bold/centered) formatting of the
# header row above it. Clear that so the new data row matches the plain
# formatting used by every other data row in the sheet.
$ws.Rows(2).ClearFormats()

# Columns A, B, C, E, F, G, H, I, J, K, R are identical for every record
# in this sheet (same market/product/origin). Copy them from the row
# that used to be row 2 (now row 3) into the freshly inserted row 2.
$repeatCols = 1,2,3,5,6,7,8,9,10,11,18
foreach ($c in $repeatCols) {
    $ws.Cells.Item(2, $c).Value2 = $ws.Cells.Item(3, $c).Value2
}

# New data specific to the inserted row 2
$ws.Cells.Item(2, 4).Value2 = 44616                       # D2 Fecha
$ws.Cells.Item(2, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS" # match date formatting used elsewhere in column D
$ws.Cells.Item(2, 12).Value2 = "Segunda"                  # L2 Calidad
$ws.Cells.Item(2, 13).Value2 = 300                        # M2 Volumen
$ws.Cells.Item(2, 14).Value2 = 16000                      # N2 Precio minimo
$ws.Cells.Item(2, 15).Value2 = 17000                      # O2 Precio maximo
$ws.Cells.Item(2, 16).Value2 = 16500                      # P2 Precio promedio ponderado
$ws.Cells.Item(2, 17).Value2 = "`$/caja 18 kilos granel"  # Q2 Unidad de comercializacion
$ws.Cells.Item(2, 19).Value2 = 917                        # S2 Precio $/Kg
$ws.Cells.Item(2, 20).Value2 = 18                         # T2 Kg / unidad
